$d = $word.ActiveDocument

# --- 1) Indent fix on the first four header paragraphs (720 -> -90 twips, i.e. 36pt -> -4.5pt) ---
for ($i = 1; $i -le 4; $i++) {
    $para = $d.Paragraphs.Item($i)
    $para.Format.LeftIndent = -4.5
}

# --- 2) Replace the {contractDate} / {contractMonth} bold placeholders with a plain "…." ellipsis ---
$rngDate = $d.Content
$foundDate = $rngDate.Find.Execute("{contractDate}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, $false)
if ($foundDate) {
    $rngDate.Delete()
    $rngDate.InsertBefore("….")
}

$rngMonth = $d.Content
$foundMonth = $rngMonth.Find.Execute("{contractMonth}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, $false)
if ($foundMonth) {
    $rngMonth.Delete()
    $rngMonth.InsertBefore("….")
}
